# Merge split runs into single runs (fixing the underlying text in the process).
#
# Slide 2, shape 2, paragraph 1:
#   "Our idea was to make FSMs to monitor/control " + "Heart Events."
#   -> "Our idea was to make FSMs to monitor/control Heart Events."
#
# Slide 5, shape 2, paragraph 5:
#   "Mode 2: ... with the " + "pacemaker.It" + " paces the heart according to diseases selected."
#   -> "Mode 2: ... with the pacemaker. It paces the heart according to diseases selected."

$p = $ppt.ActivePresentation

# --- Slide 2 -------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

$old2 = "Our idea was to make FSMs to monitor/control Heart Events."
$new2 = "Our idea was to make FSMs to monitor/control Heart Events."

$full2 = $tr2.Text
$idx2  = $full2.IndexOf($old2)
$sub2  = $tr2.Characters($idx2 + 1, $old2.Length)
$sub2.Text = $new2

# --- Slide 5 -------------------------------------------------------------
$s5  = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange

$old5 = "Mode 2: Non-blocking UART is used to interface the virtual heart program with the pacemaker.It paces the heart according to diseases selected."
$new5 = "Mode 2: Non-blocking UART is used to interface the virtual heart program with the pacemaker. It paces the heart according to diseases selected."

$full5 = $tr5.Text
$idx5  = $full5.IndexOf($old5)
$sub5  = $tr5.Characters($idx5 + 1, $old5.Length)
$sub5.Text = $new5
